$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New master-data rows for "reg-ack-template-part4" (eng / ara / fra), following
# the exact pattern already used by part1/part2/part3 above them.

# Row 122 - English
$ws.Range("A122").Value = "reg-ack-template-part4"
$ws.Range("B122").Value = "Registration Acknowledgement Template - Part 4"
$ws.Range("C122").Value = "eng"
$ws.Range("D122").Value = $true
$ws.Range("E122").Value = "superadmin"
$ws.Range("F122").Value = "now()"

# Row 123 - Arabic
$ws.Range("A123").Value = "reg-ack-template-part4"
$ws.Range("B123").Value = "نموذج شكر التسجيل"
$ws.Range("C123").Value = "ara"
$ws.Range("D123").Value = $true
$ws.Range("E123").Value = "superadmin"
$ws.Range("F123").Value = "now()"

# Row 124 - French
$ws.Range("A124").Value = "reg-ack-template-part4"
$ws.Range("B124").Value = "accusé de réception"
$ws.Range("C124").Value = "fra"
$ws.Range("D124").Value = $true
$ws.Range("E124").Value = "superadmin"
$ws.Range("F124").Value = "now()"

# Match the author's final selection/active cell after entering the new rows
$ws.Range("A125:XFD1048576").Select()
